$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 130 ("Fruta / hortaliza, semanal"
# commit): Vega Monumental Concepción - Uva, Red Globe, Primera, 2023-... entry.
# Inserting the row shifts the previous rows 130-240 down to 131-241.
$ws.Rows.Item(130).Insert()

$row = 130

$ws.Cells.Item($row, 1).Value  = 11
$ws.Cells.Item($row, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value  = "Bíobío"
$ws.Cells.Item($row, 4).Value  = 45049
$ws.Cells.Item($row, 5).Value  = 8
$ws.Cells.Item($row, 6).Value  = "Fruta"
$ws.Cells.Item($row, 7).Value  = 100109
$ws.Cells.Item($row, 8).Value  = "Uva"
$ws.Cells.Item($row, 9).Value  = 100109001
$ws.Cells.Item($row, 10).Value = "Uva"
$ws.Cells.Item($row, 11).Value = "Red Globe"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 150
$ws.Cells.Item($row, 14).Value = 9000
$ws.Cells.Item($row, 15).Value = 10000
$ws.Cells.Item($row, 16).Value = 9533
$ws.Cells.Item($row, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 530
$ws.Cells.Item($row, 20).Value = 18
